# Change the original power buy price from 0.24 to 0.2 (20 cents) for four
# daily price blocks (rows 182-277, 326-421, 470-565, 614-708 in column B).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B182:B277").Value = 0.2
$ws.Range("B326:B421").Value = 0.2
$ws.Range("B470:B565").Value = 0.2
$ws.Range("B614:B708").Value = 0.2

# Reflect the author's final cursor/scroll position recorded in the sheet view.
$ws.Range("D710").Select()
